$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Sheet is protected; unprotect to allow edits, then restore protection after
$ws.Unprotect()

# Update the "as of" date in the confidential disclaimer text (cell A10)
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-18 for illustrative purposes only and are subject to change."

# Refresh the Weight (D) and Percent Change (E) values for rows 2-7
$ws.Range("D2").Value = 0.241762274239199
$ws.Range("E2").Value = -0.003091590057882287

$ws.Range("D3").Value = 0.5040963733422339
$ws.Range("E3").Value = -0.01044386422976507

$ws.Range("D4").Value = 0.09441725471823451
$ws.Range("E4").Value = -0.003025413473174532

$ws.Range("D5").Value = 0.1031016740733343
$ws.Range("E5").Value = -0.01022329835889157

$ws.Range("D6").Value = 0.05662242362699826
$ws.Range("E6").Value = -0.007481296758104716

$ws.Range("D7").Value = 0.9999999999999999
$ws.Range("E7").Value = -0.007775443489506539

# Restore sheet protection (original allowed column/row formatting while locked)
$ws.Protect($null, $true, $true, $true, $false, $false, $true, $true, $false, $false, $false, $false, $false, $false, $false, $false)
